# Update the values in the "two-digit number / one-digit number" division
# worksheet table. Each data row of the table holds 5 expressions in
# individual cells.
#
# NOTE: this runtime's Range.Find.Execute does not stay confined to the
# range/cell it is invoked on -- it always matches the first occurrence
# in the whole document. Since several of the new values coincide with
# old values that are still pending replacement elsewhere in the table
# (e.g. "36÷3=" is both a value being replaced and a value being written
# in a different cell), a Find/Replace-based approach can clobber the
# wrong cell. Instead we set each cell's text directly via its precise
# Range (excluding the trailing end-of-cell mark), which is guaranteed
# to touch only that single cell.

$d = $word.ActiveDocument
$t = $d.Tables.Item(1)

function Set-CellText($table, $row, $col, $newText) {
    $cell = $table.Cell($row, $col)
    $full = $cell.Range
    $textRange = $d.Range($full.Start, $full.End - 1)
    $textRange.Text = $newText
}

# row -> new values for each of the 5 columns, in column order
$rowData = @{
    1  = @("19÷4=", "12÷3=", "36÷3=", "84÷8=", "94÷2=")
    5  = @("44÷9=", "86÷4=", "38÷7=", "92÷2=", "14÷4=")
    9  = @("21÷6=", "78÷9=", "20÷5=", "86÷3=", "91÷9=")
    13 = @("83÷8=", "20÷3=", "17÷5=", "55÷4=", "79÷3=")
    17 = @("85÷2=", "97÷2=", "32÷5=", "18÷9=", "19÷5=")
}

foreach ($rowIndex in ($rowData.Keys | Sort-Object)) {
    $values = $rowData[$rowIndex]
    for ($col = 1; $col -le $values.Count; $col++) {
        Set-CellText $t $rowIndex $col $values[$col - 1]
    }
}
